$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 15: Inscritos 146 -> 148
$ws.Range("E15").Value = 148

# Row 19: Inscritos 47 -> 48
$ws.Range("E19").Value = 48

# Row 25: Inscritos 16 -> 17, Pagos 7 -> 8, Inscricoes homologadas 7 -> 8
$ws.Range("E25").Value = 17
$ws.Range("F25").Value = 8
$ws.Range("H25").Value = 8

# Row 44: Inscritos 24 -> 25
$ws.Range("E44").Value = 25

# Row 47: Inscritos 50 -> 51
$ws.Range("E47").Value = 51

# Row 48: Pagos 17 -> 19, Inscricoes homologadas 17 -> 19
$ws.Range("F48").Value = 19
$ws.Range("H48").Value = 19

# Row 57: Inscritos 11 -> 12
$ws.Range("E57").Value = 12

# Row 76: Inscritos 45 -> 47, Pagos 16 -> 17, Inscricoes homologadas 16 -> 17
$ws.Range("E76").Value = 47
$ws.Range("F76").Value = 17
$ws.Range("H76").Value = 17

# Row 89: Inscritos 29 -> 30
$ws.Range("E89").Value = 30
